$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "codigo_partida"
$ws.Range("B1").Value = "codigo_catalogo"
$ws.Range("C1").Value = "nombre"
$ws.Range("D1").Value = "tipo"
$ws.Range("E1").Value = "descripcion"
$ws.Range("F1").Value = "cantidad_contenedor"
$ws.Range("G1").Value = "unidades_contenedor"

# Row 2 - FC51 / Mesa
$ws.Range("A2").Value = "FC51"
$ws.Range("B2").Value = "FC51"
$ws.Range("C2").Value = "Mesa"
$ws.Range("D2").Value = "Mobiliario y equipo de oficina"
$ws.Range("E2").Value = "test"
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2

# Row 3 - U1000 / Silla
$ws.Range("A3").Value = "U1000"
$ws.Range("B3").Value = "U1000"
$ws.Range("C3").Value = "Silla"
$ws.Range("D3").Value = "Mobiliario y equipo de oficina"
$ws.Range("E3").Value = "test"
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 50

# Row 4 - U2000 / Computadora (no descripcion)
$ws.Range("A4").Value = "U2000"
$ws.Range("B4").Value = "U2000"
$ws.Range("C4").Value = "Computadora"
$ws.Range("D4").Value = "Equipos informáticos"
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5

# Row 5 - U3000 / Laptop
$ws.Range("A5").Value = "U3000"
$ws.Range("B5").Value = "U3000"
$ws.Range("C5").Value = "Laptop"
$ws.Range("D5").Value = "Equipos informáticos"
$ws.Range("E5").Value = "test"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1

# Column widths for the new columns (closest achievable values on this
# engine's char-width grid to the authored widths of 19.28515625 / 26.140625)
$ws.Columns.Item(5).ColumnWidth = 18.5
$ws.Columns.Item(6).ColumnWidth = 25.33

# Selection state
$ws.Range("F10").Select()
